$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B11 currently holds the text "R40" (a shared string). The edit renames
# this rule's label to the text "1". Assigning a numeric-looking string
# straight to .Value would auto-convert it to a real number, so instead we
# write it as a formula that evaluates to a text "1" and then collapse the
# formula down to a plain text value in-place (Copy + PasteSpecial values),
# exactly like pasting values-only in the Excel UI. This keeps the cell's
# existing style/number format untouched and stores "1" as a shared string.
$ws.Range("B11").Formula = "=""1"""
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

